$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlinks -----------------------------------------------------------
# Row 3's hyperlinks (J3, K3) must be removed, while J2/K2's hyperlinks must
# stay untouched. In this COM host, Hyperlinks.Delete() (whether reached via
# the worksheet or a range) always removes *every* hyperlink on the sheet
# (matches real Excel's well-known quirk), and deleting an individual
# Hyperlink item is a no-op. So: stash J2/K2's link info + cell formatting,
# wipe all hyperlinks, then recreate only the ones we want to keep and
# restore their original formatting (avoids leaving behind a duplicated
# cell style record).
$keepAddr = @{}
$keepTip = @{}
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address($false, $false)
    if ($addr -eq "J2" -or $addr -eq "K2") {
        $keepAddr[$addr] = $h.Address
        $keepTip[$addr] = $h.ScreenTip
    }
}

$ws.Range("J2:K2").Copy($ws.Range("M2:N2"))

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("J2"), $keepAddr["J2"], [Type]::Missing, $keepTip["J2"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K2"), $keepAddr["K2"], [Type]::Missing, $keepTip["K2"]) | Out-Null

$ws.Range("M2:N2").Copy()
$ws.Range("J2:K2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("M2:N2").Clear()
$excel.CutCopyMode = $false

# --- Row 3 contents ---------------------------------------------------------
# Clear the contents of A3:I3 entirely (no more value/type on these cells)
$ws.Range("A3:I3").ClearContents()

# J3/K3 keep their (hyperlink) style/formatting, only the value/content and
# the hyperlink itself are removed
$ws.Range("J3:K3").ClearContents()

# --- Selection --------------------------------------------------------------
$ws.Range("B2").Select()
